$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.946.22"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "'2.667.08"
$ws.Range("E3").Value = "  +2.22%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'574.86"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "'144.40"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("E11").Value = "  +2.79%  "

$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "'3.136.77"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").Value = "'26.25"
$ws.Range("E14").Value = "  +11.81%  "

$ws.Range("D15").Value = "'60.949.28"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("D17").Value = "'2.669.82"
$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("E18").Value = "  +2.57%  "

$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").Value = "'350.92"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").Value = "'63.98"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").Value = "'8.25"
$ws.Range("E27").Value = "  +3.31%  "

$ws.Range("E28").Value = "  +9.25%  "

$ws.Range("D29").Value = "'0.0₃0809"
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +7.23%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").Value = "'163.49"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("E34").Value = "  +8.45%  "

$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("E36").Value = "  +6.68%  "

$ws.Range("D37").Value = "'340.67"
$ws.Range("E37").Value = "  +11.20%  "

$ws.Range("E38").Value = "  +2.14%  "

$ws.Range("D39").Value = "'4.10"
$ws.Range("E39").Value = "  +5.45%  "

$ws.Range("E40").Value = "  +6.69%  "

$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("E42").Value = "  +2.89%  "

$ws.Range("D44").Value = "'20.37"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").Value = "'132.99"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'20.61"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0996"
$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").Value = "'2.099.32"
$ws.Range("E51").Value = "  +3.21%  "
